# Update part numbers for 7x7 maze IR module cables.
#
# A new line item for "Molex Micro-Fit 3.0 cable 3m" (qty 1) needs to be
# added directly above the existing "Molex Micro-Fit 3.0 cable 2m" row,
# and the quantity of the 2m cable row drops from 4 to 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the existing "Molex Micro-Fit 3.0 cable 2m"
# row (row 84). This shifts that row and everything below it down by one.
$ws.Rows.Item(84).Insert() | Out-Null

# Fill in the details for the newly inserted row (new cable part).
$ws.Cells.Item(84, 1).Value2 = "Molex Micro-Fit 3.0 cable 3m"
$ws.Cells.Item(84, 2).Value2 = "IR lighting"
$ws.Cells.Item(84, 3).Value2 = 1
$ws.Cells.Item(84, 4).Value2 = "Farnell"
$ws.Cells.Item(84, 5).Value2 = 3862421

# The original "Molex Micro-Fit 3.0 cable 2m" row (now shifted to row 85)
# has its quantity reduced from 4 to 3.
$ws.Cells.Item(85, 3).Value2 = 3

# Reflect the cell the user ended up on after making the edit.
$ws.Range("B95").Select() | Out-Null
